$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing rows 2-4 down to 3-5)
$ws.Rows.Item(2).Insert()

# Insert a new row at the end (after what is now row 5) to become row 6
$ws.Rows.Item(6).Insert()

# New row 2 data
$ws.Cells.Item(2, 1).Value = "عاشق"
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 4
$ws.Cells.Item(2, 4).Value = "2022-04-17 11:36:00"
$ws.Cells.Item(2, 5).Value = "اونایی که در خط امام نیستن"
$ws.Cells.Item(2, 6).Value = "Open"
$ws.Cells.Item(2, 7).Value = "alireza"

# New row 6 data
$ws.Cells.Item(6, 1).Value = "عاشق"
$ws.Cells.Item(6, 2).Value = 5
$ws.Cells.Item(6, 3).Value = 10
$ws.Cells.Item(6, 4).Value = "2022-04-19 08:16:11"
$ws.Cells.Item(6, 5).Value = "خرابی کد تخفیف"
$ws.Cells.Item(6, 6).Value = "Answered"
$ws.Cells.Item(6, 7).Value = "mamd"
